$d = $word.ActiveDocument

# --- Locate the QUESTION1 / QUESTION2 paragraphs in the "Questions" list ---
# Paragraph.Range.Text includes the trailing paragraph mark (CR), so trim it
# before comparing.
$q1 = $null
$q2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "QUESTION1") { $q1 = $p }
    if ($t -eq "QUESTION2") { $q2 = $p }
}

# Work bottom-up so inserting after QUESTION1 doesn't disturb QUESTION2's Range.
if ($q2 -ne $null) {
    $q2.Range.InsertAfter([char]13 + "ANSWER")
    $ansPara2 = $q2.Next()
    $ansPara2.Range.ListFormat.ListLevelNumber = 2
}

if ($q1 -ne $null) {
    $q1.Range.InsertAfter([char]13 + "ANSWER")
    $ansPara1 = $q1.Next()
    $ansPara1.Range.ListFormat.ListLevelNumber = 2
}

# --- Page orientation: explicit portrait on the section(s) ---
foreach ($sec in $d.Sections) {
    $sec.PageSetup.Orientation = 0
}
